$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value2 = 3.2
$ws.Range("I3").Value2 = 2.4
$ws.Range("V3").Value2 = 12
$ws.Range("Y3").Value2 = 41
$ws.Range("AG3").Value2 = 9.5
$ws.Range("J6").Value2 = 1.11
$ws.Range("K6").Value2 = 6.5
$ws.Range("N7").Value2 = 2.1
$ws.Range("O7").Value2 = 1.7
$ws.Range("G8").Value2 = 1.33
$ws.Range("H8").Value2 = 5.5
$ws.Range("V8").Value2 = 8.5
$ws.Range("Z8").Value2 = 13
$ws.Range("AG8").Value2 = 23
$ws.Range("L9").Value2 = 1.2
$ws.Range("M9").Value2 = 4.5
$ws.Range("N9").Value2 = 1.67
$ws.Range("O9").Value2 = 2.2
$ws.Range("G10").Value2 = 2.5
$ws.Range("H10").Value2 = 3.1
$ws.Range("I10").Value2 = 3
$ws.Range("P10").Value2 = 1.62
$ws.Range("Q10").Value2 = 2.2
$ws.Range("R10").Value2 = 2.2
$ws.Range("S10").Value2 = 1.62
$ws.Range("AC10").Value2 = 81
$ws.Range("AE10").Value2 = 7
$ws.Range("AF10").Value2 = 13
$ws.Range("N11").Value2 = 2.06
$ws.Range("O11").Value2 = 1.84
$ws.Range("G20").Value2 = 2.42
$ws.Range("H20").Value2 = 2.87
$ws.Range("I20").Value2 = 3.05
$ws.Range("N20").Value2 = 2.27
$ws.Range("O20").Value2 = 1.5
$ws.Range("AA20").Value2 = 5.7
$ws.Range("T20").Value2 = 6.5
$ws.Range("U20").Value2 = 11
$ws.Range("V20").Value2 = 9.5
$ws.Range("W20").Value2 = 26
$ws.Range("X20").Value2 = 23
$ws.Range("Z20").Value2 = 6.7
$ws.Range("AE20").Value2 = 7.3
$ws.Range("AF20").Value2 = 14.5
$ws.Range("AG20").Value2 = 11.5
$ws.Range("AH20").Value2 = 40
$ws.Range("AI20").Value2 = 32
$ws.Range("AJ20").Value2 = 50
$ws.Range("AD21").Value2 = 1000
$ws.Range("I23").Value2 = 7
$ws.Range("N23").Value2 = 2.1
$ws.Range("O23").Value2 = 1.7
$ws.Range("R23").Value2 = 2.2
$ws.Range("S23").Value2 = 1.62
$ws.Range("T23").Value2 = 5.5
$ws.Range("U23").Value2 = 6
$ws.Range("V23").Value2 = 9
$ws.Range("Z23").Value2 = 8.5
$ws.Range("AC23").Value2 = 81
$ws.Range("AE23").Value2 = 15
$ws.Range("AF23").Value2 = 34
$ws.Range("AH23").Value2 = 81
$ws.Range("G27").Value2 = 4.2
$ws.Range("H27").Value2 = 3.2
$ws.Range("I27").Value2 = 1.78
$ws.Range("N27").Value2 = 1.98
$ws.Range("O27").Value2 = 1.65
$ws.Range("AA27").Value2 = 5.6
$ws.Range("AB27").Value2 = 12.5
$ws.Range("AC27").Value2 = 55
$ws.Range("AD27").Value2 = 400
$ws.Range("AE27").Value2 = 5.4
$ws.Range("AF27").Value2 = 6.7
$ws.Range("AG27").Value2 = 7
$ws.Range("AH27").Value2 = 11.75
$ws.Range("AI27").Value2 = 12.5
$ws.Range("AJ27").Value2 = 23
$ws.Range("R27").Value2 = 1.8
$ws.Range("S27").Value2 = 1.9
$ws.Range("T27").Value2 = 9.5
$ws.Range("U27").Value2 = 19.5
$ws.Range("V27").Value2 = 11.5
$ws.Range("W27").Value2 = 55
$ws.Range("X27").Value2 = 32
$ws.Range("Y27").Value2 = 35
$ws.Range("Z27").Value2 = 8.5
$ws.Range("N28").Value2 = 1.75
$ws.Range("O28").Value2 = 2.05
$ws.Range("P28").Value2 = 1.33
$ws.Range("Q28").Value2 = 3.25
$ws.Range("R28").Value2 = 2.2
$ws.Range("S28").Value2 = 1.62
$ws.Range("T28").Value2 = 6.5
$ws.Range("U28").Value2 = 6
$ws.Range("AI28").Value2 = 67
$ws.Range("H30").Value2 = 5.5
$ws.Range("I30").Value2 = 7.3
$ws.Range("N30").Value2 = 1.3
$ws.Range("O30").Value2 = 2.95
$ws.Range("R30").Value2 = 1.55
$ws.Range("S30").Value2 = 2.15
$ws.Range("T30").Value2 = 12.5
$ws.Range("U30").Value2 = 9.25
$ws.Range("W30").Value2 = 10
$ws.Range("Y30").Value2 = 18.5
$ws.Range("J33").Value2 = 1.08
$ws.Range("K33").Value2 = 8
$ws.Range("L33").Value2 = 1.44
$ws.Range("M33").Value2 = 2.63
$ws.Range("G34").Value2 = 2.15
$ws.Range("I34").Value2 = 3.6
$ws.Range("U34").Value2 = 9
$ws.Range("AE34").Value2 = 8
$ws.Range("J35").Value2 = 1.07
$ws.Range("K35").Value2 = 9
$ws.Range("N35").Value2 = 2.25
$ws.Range("O35").Value2 = 1.62
$ws.Range("I36").Value2 = 3.95
$ws.Range("K36").Value2 = 8.75
$ws.Range("L36").Value2 = 1.19
$ws.Range("M36").Value2 = 4.15
$ws.Range("N36").Value2 = 1.6
$ws.Range("O36").Value2 = 2.22
$ws.Range("P36").Value2 = 1.31
$ws.Range("Q36").Value2 = 3.15
$ws.Range("R36").Value2 = 1.55
$ws.Range("S36").Value2 = 2.3
$ws.Range("AA36").Value2 = 7.5
$ws.Range("AB36").Value2 = 12.5
$ws.Range("AC36").Value2 = 40
$ws.Range("AD36").Value2 = 250
$ws.Range("AE36").Value2 = 15.5
$ws.Range("AF36").Value2 = 26
$ws.Range("W36").Value2 = 15.5
$ws.Range("Y36").Value2 = 19.5
$ws.Range("Z36").Value2 = 8.75
$ws.Range("AJ36").Value2 = 30
$ws.Range("G37").Value2 = 2.67
$ws.Range("H37").Value2 = 3.85
$ws.Range("I37").Value2 = 2.25
$ws.Range("M37").Value2 = 4.9
$ws.Range("N37").Value2 = 1.44
$ws.Range("O37").Value2 = 2.57
$ws.Range("Q37").Value2 = 3.55
$ws.Range("S37").Value2 = 2.67
$ws.Range("T37").Value2 = 14.5
$ws.Range("U37").Value2 = 18
$ws.Range("W37").Value2 = 32
$ws.Range("X37").Value2 = 18.5
$ws.Range("AD37").Value2 = 150
$ws.Range("AH37").Value2 = 24
$ws.Range("G42").Value2 = 2.32
$ws.Range("H42").Value2 = 3.15
$ws.Range("I42").Value2 = 3
$ws.Range("N42").Value2 = 2.37
$ws.Range("O42").Value2 = 1.52
$ws.Range("R42").Value2 = 2.1
$ws.Range("S42").Value2 = 1.65
$ws.Range("T42").Value2 = 6.1
$ws.Range("U42").Value2 = 9.75
$ws.Range("V42").Value2 = 9.75
$ws.Range("W42").Value2 = 23
$ws.Range("X42").Value2 = 23
$ws.Range("AA42").Value2 = 6.2
$ws.Range("AB42").Value2 = 18.5
$ws.Range("AC42").Value2 = 120
$ws.Range("AE42").Value2 = 7.2
$ws.Range("AF42").Value2 = 13.5
$ws.Range("AG42").Value2 = 11.5
$ws.Range("AH42").Value2 = 37
$ws.Range("AI42").Value2 = 32
$ws.Range("AJ42").Value2 = 50
$ws.Range("G47").Value2 = 1.37
$ws.Range("H47").Value2 = 4.2
$ws.Range("I47").Value2 = 9.25
$ws.Range("J47").Value2 = 1.06
$ws.Range("K47").Value2 = 7.5
$ws.Range("L47").Value2 = 1.27
$ws.Range("P47").Value2 = 1.4
$ws.Range("Q47").Value2 = 2.75
$ws.Range("T47").Value2 = 6
$ws.Range("U47").Value2 = 5.9
$ws.Range("Z47").Value2 = 7.5
$ws.Range("AD47").Value2 = 900
$ws.Range("AE47").Value2 = 22
$ws.Range("L56").Value2 = 1.23
$ws.Range("M56").Value2 = 3.4
$ws.Range("N56").Value2 = 1.8
$ws.Range("O56").Value2 = 2
$ws.Range("J58").Value2 = 1.07
$ws.Range("K58").Value2 = 9
$ws.Range("N58").Value2 = 2.1
$ws.Range("O58").Value2 = 1.7
$ws.Range("G62").Value2 = 2.45
$ws.Range("H62").Value2 = 2.9
$ws.Range("I62").Value2 = 2.75
$ws.Range("R62").Value2 = 2.2
$ws.Range("S62").Value2 = 1.62
$ws.Range("T62").Value2 = 6.5
$ws.Range("U62").Value2 = 11
$ws.Range("V62").Value2 = 11
$ws.Range("W62").Value2 = 26
$ws.Range("X62").Value2 = 26
$ws.Range("AA62").Value2 = 6
$ws.Range("AE62").Value2 = 6.5
$ws.Range("AF62").Value2 = 12
$ws.Range("AG62").Value2 = 12
$ws.Range("AH62").Value2 = 29
$ws.Range("AI62").Value2 = 29
$ws.Range("AJ62").Value2 = 41
$ws.Range("H64").Value2 = 3.4
$ws.Range("J64").Value2 = 1.01
$ws.Range("K64").Value2 = 15
$ws.Range("AA64").Value2 = 7
$ws.Range("J66").Value2 = 1.04
$ws.Range("K66").Value2 = 9
$ws.Range("L66").Value2 = 1.22
$ws.Range("M66").Value2 = 4
$ws.Range("N66").Value2 = 1.7
$ws.Range("O66").Value2 = 2.1
$ws.Range("H69").Value2 = 3.95
$ws.Range("I69").Value2 = 5.7
$ws.Range("T69").Value2 = 7.5
$ws.Range("U69").Value2 = 7.5
$ws.Range("V69").Value2 = 8
$ws.Range("W69").Value2 = 10.75
$ws.Range("X69").Value2 = 11.5
$ws.Range("N73").Value2 = 1.67
$ws.Range("O73").Value2 = 2.15
$ws.Range("L74").Value2 = 1.18
$ws.Range("M74").Value2 = 4.5
$ws.Range("N74").Value2 = 1.6
$ws.Range("O74").Value2 = 2.3
$ws.Range("M80").Value2 = 3.2
$ws.Range("N80").Value2 = 1.9
$ws.Range("O80").Value2 = 1.8
$ws.Range("P80").Value2 = 1.42
$ws.Range("Q80").Value2 = 2.65
$ws.Range("R80").Value2 = 1.72
$ws.Range("S80").Value2 = 2
$ws.Range("AB80").Value2 = 13.5
$ws.Range("AC80").Value2 = 60
$ws.Range("AD80").Value2 = 450
$ws.Range("AE80").Value2 = 9.75
